$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The source data for rows 83 and 84 (BW Linz vs LASK, and A. Klagenfurt vs
# Sturm Graz) got re-sorted: the match details (columns F:V — teams, score,
# odds, timestamps, url) that used to sit on row 83 now belong on row 84 and
# vice versa. Column A (the running "Indice") and B:D stay put since both
# rows already share the same country/tournament/season. Swap the F:V blocks
# between the two rows.
# ---------------------------------------------------------------------------
$row83 = $ws.Range("F83:V83").Value2
$row84 = $ws.Range("F84:V84").Value2

$ws.Range("F83:V83").Value2 = $row84
$ws.Range("F84:V84").Value2 = $row83

# ---------------------------------------------------------------------------
# A brand-new match (Austria Vienna vs Salzburg) was appended as row 85.
# Copy the formatting (number format / style) of the styled columns (A, E)
# from row 84 so the new row matches the existing look, then fill in values.
# ---------------------------------------------------------------------------
$ws.Range("A84").Copy($ws.Range("A85"))
$ws.Range("E84").Copy($ws.Range("E85"))

$ws.Range("A85").Value = 84
$ws.Range("B85").Value = "austria"
$ws.Range("C85").Value = "bundesliga"
$ws.Range("D85").Value = "2023-2024"
$ws.Range("E85").Value = 45242.70833333334
$ws.Range("F85").Value = "Austria Vienna"
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = "Salzburg"
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 4.84
$ws.Range("K85").Value = "04/11/2023 17:42"
$ws.Range("L85").Value = 4.26
$ws.Range("M85").Value = "12/11/2023 16:58"
$ws.Range("N85").Value = 4.39
$ws.Range("O85").Value = "04/11/2023 17:42"
$ws.Range("P85").Value = 3.85
$ws.Range("Q85").Value = "12/11/2023 16:57"
$ws.Range("R85").Value = 1.65
$ws.Range("S85").Value = "04/11/2023 17:42"
$ws.Range("T85").Value = 1.87
$ws.Range("U85").Value = "12/11/2023 16:57"
$ws.Range("V85").Value = "https://www.betexplorer.com/football/austria/bundesliga/austria-vienna-salzburg/jsiV8ASI/"
